# Update "想去人数" (F column) counts on three worksheets to reflect the
# newly generated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    6  = 664
    8  = 2731
    10 = 6367
    11 = 2415
    13 = 26
    15 = 2872
    18 = 6840
    19 = 252
    21 = 188
    24 = 7768
    27 = 251
    31 = 112
    32 = 51
    34 = 11
    35 = 31
    37 = 66
    38 = 2578
    41 = 30
    44 = 601
    45 = 3615
    46 = 143
    47 = 1163
    48 = 104
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    2  = 33
    5  = 236
    7  = 107
    15 = 164
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Range("F$row").Value = $sheet2Updates[$row]
}

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    4  = 664
    6  = 33
    7  = 2731
    9  = 236
    10 = 236
    11 = 6367
    12 = 2415
    13 = 107
    14 = 26
    16 = 2872
    21 = 6840
    22 = 252
    24 = 188
    27 = 7768
    29 = 251
    33 = 112
    34 = 11
    35 = 31
    37 = 66
    38 = 2578
    41 = 30
    44 = 601
    45 = 164
    46 = 3615
    47 = 143
    49 = 1163
    50 = 104
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
